# Consolidated error messages V2.3 - update IDA-MLC error message rows
# (see commit: "Updated error message for IDA")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IDA")

# Row 14 (scenario #13): drop the "; Failure in Decryption" suffix from the
# "Could not process request/Unknown error..." scenario text.
$ws.Range("B14").Value = "Could not process request/Unknown error; Invalid Auth Request"

# Row 34 (scenario #40): shorten the message text and highlight the row in
# yellow to flag the update (B34 keeps its original scenario text, C34 gets
# the new, shorter message).
$ws.Range("B34").Value = "Invalid encryption of session key/request"
$ws.Range("C34").Value = [char]0x201C + "Unable to decrypt Request." + [char]0x201D
$ws.Range("B34:F34").Interior.Color = 65535

# Selection moved back up to the top of the sheet (was scrolled to A58/J9).
$ws.Range("J7").Select() | Out-Null
